# Add a new "2022-Q1" worksheet before the "总计" (Total) sheet,
# populate it with the fund-holding detail rows, and update the "总计"
# summary sheet with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned immediately before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. Header row for the new sheet
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160
$newSheet.Range("B1:H1").Borders.LineStyle = 1

$q1Data = @(
    @(0, '007119', '睿远成长价值混合A', '280.90', '92.63', '3.72', '10.4495', 9),
    @(1, '000751', '嘉实新兴产业股票', '81.04', '90.66', '8.12', '6.5804', 1),
    @(2, '010186', '嘉实核心成长混合A', '69.42', '90.46', '8.64', '5.9979', 1),
    @(3, '009795', '嘉实远见精选两年持有期混合', '63.75', '89.19', '7.95', '5.0681', 1),
    @(4, '002593', '富国美丽中国混合', '55.22', '87.05', '7.31', '4.0366', 2),
    @(5, '100026', '富国天合稳健混合', '53.82', '80.56', '7.11', '3.8266', 1),
    @(6, '000595', '嘉实泰和混合', '42.72', '90.51', '7.94', '3.3920', 1),
    @(7, '009137', '嘉实瑞和两年持有期混合', '26.52', '87.46', '7.77', '2.0606', 1),
    @(8, '070002', '嘉实增长混合', '27.30', '71.00', '6.98', '1.9055', 1),
    @(9, '007120', '睿远成长价值混合C', '29.98', '92.63', '3.72', '1.1153', 9),
    @(10, '519035', '富国天博创新混合', '25.63', '93.05', '4.05', '1.0380', 5),
    @(11, '166025', '中欧远见两年定期开放混合A', '48.80', '48.74', '2.04', '0.9955', 10),
    @(12, '360006', '光大保德信新增长混合', '21.71', '88.07', '3.98', '0.8641', 8),
    @(13, '000513', '富国高端制造行业股票', '9.25', '93.54', '4.49', '0.4153', 4),
    @(14, '070022', '嘉实领先成长混合', '5.27', '88.31', '7.48', '0.3942', 2),
    @(15, '010187', '嘉实核心成长混合C', '4.50', '90.46', '8.64', '0.3888', 1),
    @(16, '001740', '光大保德信中国制造2025灵活配置混合', '11.43', '86.23', '3.33', '0.3806', 6),
    @(17, '001759', '嘉实成长增强灵活配置混合', '4.59', '90.80', '7.48', '0.3433', 2),
    @(18, '160727', '嘉实创业板两年定期开放混合', '10.97', '64.85', '3.07', '0.3368', 10),
    @(19, '007016', '富国睿泽回报混合', '6.81', '73.71', '4.76', '0.3242', 7),
    @(20, '001036', '嘉实企业变革股票', '4.03', '90.85', '7.00', '0.2821', 1),
    @(21, '011921', '富国均衡成长三年持有期混合A', '7.52', '93.77', '3.14', '0.2361', 8),
    @(22, '001040', '新华策略精选股票', '6.15', '93.72', '3.09', '0.1900', 8),
    @(23, '008138', '富国龙头优势混合', '5.02', '93.93', '3.45', '0.1732', 9),
    @(24, '008313', '光大保德信研究精选混合', '2.64', '88.98', '4.33', '0.1143', 7),
    @(25, '003292', '嘉实优势成长灵活配置混合', '1.04', '92.17', '7.90', '0.0822', 1),
    @(26, '002311', '创金合信中证500指数增强A', '5.72', '92.95', '1.14', '0.0652', 3),
    @(27, '007101', '中欧远见两年定期开放混合C', '2.87', '48.74', '2.04', '0.0585', 10),
    @(28, '168101', '九泰锐智事件驱动混合（LOF）', '0.73', '89.22', '4.68', '0.0342', 10),
    @(29, '005259', '建信龙头企业股票', '1.04', '83.45', '2.98', '0.0310', 9),
    @(30, '002316', '创金合信中证500指数增强C', '2.60', '92.95', '1.14', '0.0296', 3),
    @(31, '011922', '富国均衡成长三年持有期混合C', '0.50', '93.77', '3.14', '0.0157', 8),
    @(32, '530019', '建信社会责任混合', '0.19', '77.77', '2.93', '0.0056', 9),
    @(33, '001781', '建信现代服务业股票', '0.15', '82.90', '3.31', '0.0050', 8),
    @(34, '159932', '大成中证500深市ETF', '0.42', '97.30', '0.92', '0.0039', 9),
    @(35, '010066', '方正富邦中证500指数增强A', '0.02', '91.02', '0.62', '0.0001', 6),
    @(36, '010067', '方正富邦中证500指数增强C', '0.01', '91.02', '0.62', '0.0001', 6)
)

# ---------------------------------------------------------------------------
# 3. Fill in the 37 fund-holding detail rows (row 2 .. row 38)
#    Column B (fund code) and columns D, E, F, G are stored as text (to
#    match the source data, which keeps these figures as plain text rather
#    than numeric values - this also preserves leading zeros in fund codes).
# ---------------------------------------------------------------------------
$newSheet.Range("B2:B38").NumberFormat = "@"
$newSheet.Range("D2:G38").NumberFormat = "@"

for ($i = 0; $i -lt $q1Data.Length; $i++) {
    $r = $i + 2
    $row = $q1Data[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 1).Font.Bold = $true
    $newSheet.Cells.Item($r, 1).HorizontalAlignment = -4108
    $newSheet.Cells.Item($r, 1).VerticalAlignment = -4160
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}


# ---------------------------------------------------------------------------
# 4. Update the "总计" (Total) sheet: insert a new leading data row for
#    2022-Q1 and push the previously existing rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$oldRows = @()
for ($r = 2; $r -le 6; $r++) {
    $oldRows += ,@(
        $totalSheet.Cells.Item($r, 1).Value2,
        $totalSheet.Cells.Item($r, 2).Value2,
        $totalSheet.Cells.Item($r, 3).Value2,
        $totalSheet.Cells.Item($r, 4).Value2
    )
}

# Copy the formatting of the last existing row down into the new row 7,
# since that row did not exist before and needs the bold index-number style.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

for ($i = 0; $i -lt $oldRows.Length; $i++) {
    $r = $i + 3
    $totalSheet.Cells.Item($r, 1).Value = $oldRows[$i][0]
    $totalSheet.Cells.Item($r, 2).Value = $oldRows[$i][1]
    $totalSheet.Cells.Item($r, 3).Value = $oldRows[$i][2]
    $totalSheet.Cells.Item($r, 4).Value = $oldRows[$i][3]
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 37
$totalSheet.Cells.Item(2, 4).Value = 51.24
